$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-BookRow($row, $titulo, $autor, $anio, $editorial, $genero, $codigo) {
    # F (Genero) and H (Codigo) already carry the row style; only
    # A (Titulo), B (Autor), D (Anio) and E (Editorial) still need it,
    # so copy it over from the existing styled F cell in this row.
    $ws.Cells.Item($row, 6).Copy() | Out-Null
    foreach ($col in @(1, 2, 4, 5)) {
        $ws.Cells.Item($row, $col).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    }
    $excel.CutCopyMode = 0

    $ws.Cells.Item($row, 1).Value = $titulo
    $ws.Cells.Item($row, 2).Value = $autor
    $ws.Cells.Item($row, 4).Value = $anio
    $ws.Cells.Item($row, 5).Value = $editorial
    $ws.Cells.Item($row, 6).Value = $genero
    $ws.Cells.Item($row, 8).Value = $codigo
}

Set-BookRow 27 'La ceremonia del adiós' 'Simeone de Beauvoir' 1982 'EDHASA' 'BIOGRAFIA' 'BI C 26'
Set-BookRow 28 'La carrera del Doris Hart' 'Vicki Baum' 1956 'Juventud' 'BIOGRAFIA' 'BI C 27'
Set-BookRow 29 'Diario de mi vida' 'Maria Bashkirtseff' 1948 'Espasa Calpe' 'BIOGRAFIA' 'BI C 28'
Set-BookRow 30 'Mi diario' 'Leon Bloy' 1947 'Mundo moderno' 'BIOGRAFIA' 'BI C 29'
Set-BookRow 31 'Memorias de Dolly Morton' 'Charles Carrington' 1970 'Edasa' 'BIOGRAFIA' 'BI C 30'
Set-BookRow 32 'El marqués de Sade' 'Simeone de Beauvoir' 1956 'Leviatán' 'BIOGRAFIA' 'BI C 31'
Set-BookRow 33 'Realidad y fantasia en Naguib Mahfuz' 'Mercedes del Amo' 1991 'Universidad de Granada' 'BIOGRAFIA' 'BI C 32'
Set-BookRow 34 'Mi corazón al desnudo' 'Charles Baudelaire' 1947 'Apolo' 'BIOGRAFIA' 'BI C 33'
Set-BookRow 35 'Memorias del condado de Hecate' 'Edmund Wilson' 1989 'Versal travesías' 'BIOGRAFIA' 'BI C 34'
Set-BookRow 36 'Final de cuentas' 'Simeone de Beauvoir' 1972 'Sudamericana' 'BIOGRAFIA' 'BI C 35'
Set-BookRow 37 'La plenitud de la vida' 'Simeone de Beauvoir' 1962 'Sudamericana' 'BIOGRAFIA' 'BI C 36'
Set-BookRow 38 'La vida del doctor Samuel Johnson' 'James Boswell' 1998 'Espasa Calpe' 'BIOGRAFIA' 'BI C 37'
Set-BookRow 39 'Kafka' 'Max Brod' 1951 'Emecé' 'BIOGRAFIA' 'BI C 38'
Set-BookRow 40 'La fuerza de las cosas' 'Simeone de Beauvoir' 1961 'Sudamericana' 'BIOGRAFIA' 'BI C 39'
Set-BookRow 41 'Retablo de mis recuerdos' 'Maurice Baring' 1943 'Lauro' 'BIOGRAFIA' 'BI C 40'
Set-BookRow 42 'Juaán Ginés de Sepúlveda' 'José Manuel Rodríguez Peregrina' 1993 'Universidad de Granada' 'BIOGRAFIA' 'BI C 41'
Set-BookRow 43 'Historia de una amistad' 'Jean Paul Sartre' 1965 'Nagelkop' 'BIOGRAFIA' 'BI C 42'
Set-BookRow 44 'Memorias del señor de Schnabelewopski' 'Enrique Heine' 1956 'Insula' 'BIOGRAFIA' 'BI C 43'
Set-BookRow 45 'Voltaire par lui-meme' 'René Pomeau' 1955 'Ecrivains de Toujours' 'BIOGRAFIA' 'BI C 44'
